# Correction type pour génération à partir fsh ea4a6f04ed193a83290686b2f69a3f9cd2e7f4ad
#
# - Sets the "Name" property value (row 4, column B) on the Metadata sheet,
#   which was previously left blank.
# - Refreshes the "Date" property value (row 8, column B) to reflect the
#   regeneration timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

$ws.Range("B4").Value = "MotiffinactiviteVs"
$ws.Range("B8").Value = "2025-07-18T06:40:38+00:00"
